$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "航天发展"
$ws.Range("C2").Value = "航天发展"
$ws.Range("A3").Value = "乾照光电"
$ws.Range("C3").Value = "雷科防务"
$ws.Range("A4").Value = "通宇通讯"
$ws.Range("B4").Value = "雷科防务"
$ws.Range("C4").Value = "恩捷股份"
$ws.Range("A5").Value = "雷科防务"
$ws.Range("B5").Value = "通宇通讯"
$ws.Range("C5").Value = "平潭发展"
$ws.Range("A6").Value = "中国电影"
$ws.Range("C6").Value = "国晟科技"
$ws.Range("A7").Value = "峨眉山Ａ"
$ws.Range("B7").Value = "襄阳轴承"
$ws.Range("C7").Value = "襄阳轴承"
$ws.Range("A8").Value = "平潭发展"
$ws.Range("B8").Value = "中国电影"
$ws.Range("C8").Value = "实达集团"
$ws.Range("A9").Value = "襄阳轴承"
$ws.Range("B9").Value = "平潭发展"
$ws.Range("C9").Value = "中国电影"
$ws.Range("A10").Value = "国晟科技"
$ws.Range("B10").Value = "东方精工"
$ws.Range("C10").Value = "乾照光电"
$ws.Range("A11").Value = "海欣食品"
$ws.Range("B11").Value = "银河电子"
$ws.Range("C11").Value = "峨眉山A"
$ws.Range("A12").Value = "四川金顶"
$ws.Range("B12").Value = "顺灏股份"
$ws.Range("A13").Value = "顺灏股份"
$ws.Range("B13").Value = "华电新能"
$ws.Range("C13").Value = "银河电子"
$ws.Range("A14").Value = "银河电子"
$ws.Range("B14").Value = "国晟科技"
$ws.Range("C14").Value = "航天动力"
$ws.Range("A15").Value = "东方精工"
$ws.Range("B15").Value = "海王生物"
$ws.Range("C15").Value = "海王生物"
$ws.Range("B16").Value = "普天科技"
$ws.Range("C16").Value = "海欣食品"
$ws.Range("A17").Value = "华电新能"
$ws.Range("B17").Value = "四川金顶"
$ws.Range("C17").Value = "四川金顶"
$ws.Range("A18").Value = "实达集团"
$ws.Range("B18").Value = "海欣食品"
$ws.Range("C18").Value = "合富中国"
$ws.Range("A19").Value = "航天动力"
$ws.Range("B19").Value = "航天机电"
$ws.Range("C19").Value = "特发信息"
$ws.Range("B20").Value = "通润装备"
$ws.Range("C20").Value = "茂业商业"
$ws.Range("A21").Value = "赛微电子"
$ws.Range("B21").Value = "航天动力"
$ws.Range("C21").Value = "东方精工"
